$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.314.75"
$ws.Range("E2").Value = "  -1.12%  "
$ws.Range("D3").Value = "3.320.97"
$ws.Range("E3").Value = "  +1.72%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'185.88"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.30%  "
$ws.Range("D6").Value = "'577.71"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.91%  "
$ws.Range("E7").Value = "  +1.05%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -0.21%  "
$ws.Range("D10").Value = "'6.66"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.71%  "
$ws.Range("E11").Value = "  -0.22%  "
$ws.Range("D12").Value = "3.890.77"
$ws.Range("E12").Value = "  +1.46%  "
$ws.Range("E13").Value = "  -0.69%  "
$ws.Range("D14").Value = "'27.39"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.51%  "
$ws.Range("D15").Value = "67.485.38"
$ws.Range("E15").Value = "  -0.85%  "
$ws.Range("E16").Value = "  -0.29%  "
$ws.Range("D17").Value = "3.317.64"
$ws.Range("E17").Value = "  +0.89%  "
$ws.Range("D18").Value = "'442.89"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +6.42%  "
$ws.Range("D19").Value = "'5.67"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.62%  "
$ws.Range("D20").Value = "'13.57"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +2.05%  "
$ws.Range("E21").Value = "  +2.57%  "
$ws.Range("D22").Value = "'74.23"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +3.90%  "
$ws.Range("D23").Value = "'0.998"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("D24").Value = "3.464.60"
$ws.Range("E24").Value = "  +1.69%  "
$ws.Range("D25").Value = "'0.514"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.24%  "
$ws.Range("E26").Value = "  +1.42%  "
$ws.Range("E27").Value = "  +1.42%  "
$ws.Range("D28").Value = "'9.06"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -3.95%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("E30").Value = "  +1.23%  "
$ws.Range("D31").Value = "'22.93"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.12%  "
$ws.Range("D32").Value = "'5.33"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.45%  "
$ws.Range("D33").Value = "'0.998"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.04%  "
$ws.Range("D34").Value = "'6.84"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.47%  "
$ws.Range("E35").Value = "  -0.37%  "
$ws.Range("E36").Value = "  +4.98%  "
$ws.Range("D37").Value = "'162.26"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.34%  "
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").Value = "'1.85"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.07%  "
$ws.Range("B39").Value = "EnergySwap"
$ws.Range("C39").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D39").Value = "'27.24"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.76%  "
$ws.Range("D40").Value = "2.789.87"
$ws.Range("E40").Value = "  +5.75%  "
$ws.Range("D41").Value = "'0.790"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.87%  "
$ws.Range("E42").Value = "  +0.26%  "
$ws.Range("D44").Value = "'40.27"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.32%  "
$ws.Range("E45").Value = "  -0.63%  "
$ws.Range("D46").Value = "'24.74"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.59%  "
$ws.Range("E47").Value = "  -1.39%  "
$ws.Range("D48").Value = "'326.73"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -3.07%  "
$ws.Range("E49").Value = "  -0.21%  "
$ws.Range("D50").Value = "'0.990"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.27%  "
$ws.Range("D51").Value = "'31.11"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.51%  "
